$d = $word.ActiveDocument

# The document currently has 4 paragraphs:
#   1) "Team " + "ChromeGame" (title, centered, Times New Roman 36, underline)
#   2) "Team Members: Mohit Veligenti and Evgeny Panferov"
#   3) "Project Title and Description: ChromeGame: ... Chrome Dinosaur version."
#   4) empty paragraph holding the _GoBack bookmark
#
# Target layout:
#   1) "Team Chrome Game " (plain, no formatting)
#   2) "Apple Chrome Game" (plain, no formatting)
#   3) "By Evgeny Panferov and Mohit Veligenti" (plain, with spell-check
#      proof markers around "Evgeny" and "Veligenti")
#   4) unchanged bookmark paragraph

# Remove paragraphs 2 and 3 ("Team Members..." and "Project Title...")
# entirely, leaving paragraph 1 and the trailing bookmark paragraph.
$d.Paragraphs(2).Range.Delete()
$d.Paragraphs(2).Range.Delete()

# Replace paragraph 1's whole range (text + its own formatted paragraph
# mark) with fresh, unformatted OOXML describing the three new
# paragraphs. Using raw XML insertion avoids carrying over the old
# rPr/pPr (font, size, underline, centering) from the title paragraph.
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newXml = "<w:p $wns><w:r><w:t xml:space='preserve'>Team Chrome Game </w:t></w:r></w:p>" +
          "<w:p $wns><w:r><w:t>Apple Chrome Game</w:t></w:r></w:p>" +
          "<w:p $wns>" +
            "<w:r><w:t xml:space='preserve'>By </w:t></w:r>" +
            "<w:proofErr w:type='spellStart'/>" +
            "<w:r><w:t>Evgeny</w:t></w:r>" +
            "<w:proofErr w:type='spellEnd'/>" +
            "<w:r><w:t xml:space='preserve'> Panferov and Mohit </w:t></w:r>" +
            "<w:proofErr w:type='spellStart'/>" +
            "<w:r><w:t>Veligenti</w:t></w:r>" +
            "<w:proofErr w:type='spellEnd'/>" +
          "</w:p>"

$null = $d.Paragraphs(1).Range.InsertXML($newXml)
